$d = $word.ActiveDocument

# Step 1: transform paragraph 4 ("Default" style, empty placeholder) into the
# "Heading" styled paragraph with text "Description and Relevance" and a
# Bullets character style applied to the paragraph mark.
$p4 = $d.Paragraphs.Item(4)
$xml1 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Heading"/><w:rPr><w:rStyle w:val="Bullets"/></w:rPr></w:pPr><w:r><w:t>Description and Relevance</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$p4.Range.InsertXML($xml1)

# Step 2: replace the old instructions paragraphs (now paragraphs 5-13) with
# the new report body content.
$p5 = $d.Paragraphs.Item(5)
$p13 = $d.Paragraphs.Item(13)
$r2 = $d.Range($p5.Range.Start, $p13.Range.End)
$xml2 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>Bus routes are planned to take into consideration traffic patterns, patterns of movement etc. Over a period, these assumptions change. Traffic in one part of the city may increase while traffic in another part of the city may decrease. Also, the distribution of traffic throughout the day may change. All this implies that bus schedules need constant updating to reflect the current conditions.</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">One way to approach this problem is to </w:t></w:r><w:r><w:t xml:space="preserve">consider every one-hour window and percentage the ratio of times buses </w:t></w:r><w:r><w:t>has</w:t></w:r><w:r><w:t xml:space="preserve"> been late at every stop.</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">The routes and </w:t></w:r><w:r><w:t>one-hour</w:t></w:r><w:r><w:t xml:space="preserve"> windows where we have the highest percentage of buses reporting late will be the prime candidates for replanning.</w:t></w:r></w:p><w:p><w:r><w:t>This dataset allows us to compute this information.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Heading"/></w:pPr><w:r><w:t>Novelty</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="BodyText"/></w:pPr><w:r><w:t>This problem is different from the four other problems in the assignment.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="BodyText"/></w:pPr><w:r><w:t xml:space="preserve">In the first problem, we take spatial rectangle based on the latitude and longitude. </w:t></w:r><w:r><w:t>Then for every hour, within this rectangle, we calculate how often congestion is reported.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="BodyText"/></w:pPr><w:r><w:t xml:space="preserve">In the second problem, we want to calculate the </w:t></w:r><w:r><w:t>timetable</w:t></w:r><w:r><w:t xml:space="preserve"> of a physical vehicle, as a pair of lineID and stationID, for a given day.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="BodyText"/></w:pPr><w:r><w:t>In the third problem, we compute the station that has the highest number of buses stopping at it so that we can use that station for best reach of advertisements.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="BodyText"/></w:pPr><w:r><w:t>In the fourth problem, we compute the distance travelled by individual vehicles so that we can find out when to send them for service.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="BodyText"/></w:pPr><w:r><w:t>This new problem is different from all others because we compute how many times a route has reported delays at stations, and then aggregate them by route and hour. This is not done for any of the given problems.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="BodyText"/></w:pPr></w:p><w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r2.InsertXML($xml2)

Write-Output ("Final paragraph count=" + $d.Paragraphs.Count)
